# Update Sheet1 data: replace the two personnel records with new data
# and strip the "imported" fancy formatting (hyperlink-like style on A2/A3,
# thick-bottom bordered header rows) back to a plain worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2 ----
$ws.Range("A2").Value = "ANDI SURANA"
$ws.Range("B2").Value = 12312216
$ws.Range("C2").Value = "KOMPOL"
$ws.Range("D2").Value = "ADMINISTRASI"
$ws.Range("E2").Value = "PEKALONGAN"
$ws.Range("F2").Value = "1999-02-19"
$ws.Range("G2").Value = "ISLAM"
$ws.Range("H2").Value = "JAWA"
$ws.Range("I2").Value = "2020-02-01"
$ws.Range("J2").Value = 2

# ---- Row 3 ----
$ws.Range("A3").Value = "MURNIATI"
$ws.Range("B3").Value = 12121768
$ws.Range("C3").Value = "KOMPOL"
$ws.Range("D3").Value = "ADMINISTRASI"
$ws.Range("E3").Value = "SEMARANG"
$ws.Range("F3").Value = "1999-02-10"
$ws.Range("G3").Value = "KRISTEN"
$ws.Range("H3").Value = "BUGIS"
$ws.Range("I3").Value = "2020-02-01"
$ws.Range("J3").Value = 3

# ---- Strip leftover "imported" look: hyperlink-like style on A2/A3,
# thick bottom borders/row heights on rows 1-2 ----
$ws.Range("A1:J3").ClearFormats()

$ws.Range("B1:B3,J1:J3").NumberFormat = "0"
$ws.Range("F1:F3,I1:I3").NumberFormat = "yyyy-mm-dd;@"

$ws.Rows.Item(1).RowHeight = 14.5
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 14.5

# ---- Move selection like the saved file (B3 single cell) ----
$ws.Range("B3").Select()
